# Update "想去人数" (people interested) counts in column F across sheets.
# Sheet 1 = 展览 (Exhibition), Sheet 2 = 演出 (Performance),
# Sheet 3 = 本地生活 (Local life, unchanged), Sheet 4 = 全部类型 (All types).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item(1)   # 展览
$wsShow = $wb.Worksheets.Item(2)   # 演出
$wsAll  = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 ---
$wsExpo.Cells.Item(2, 6).Value = 4738
$wsExpo.Cells.Item(3, 6).Value = 1890
$wsExpo.Cells.Item(6, 6).Value = 3192
$wsExpo.Cells.Item(7, 6).Value = 586
$wsExpo.Cells.Item(8, 6).Value = 606
$wsExpo.Cells.Item(9, 6).Value = 293
$wsExpo.Cells.Item(10, 6).Value = 659
$wsExpo.Cells.Item(11, 6).Value = 555
$wsExpo.Cells.Item(12, 6).Value = 561
$wsExpo.Cells.Item(16, 6).Value = 1393
$wsExpo.Cells.Item(18, 6).Value = 1655
$wsExpo.Cells.Item(32, 6).Value = 4048
$wsExpo.Cells.Item(33, 6).Value = 15
$wsExpo.Cells.Item(36, 6).Value = 1716
$wsExpo.Cells.Item(37, 6).Value = 61
$wsExpo.Cells.Item(38, 6).Value = 1915

# --- Sheet 2: 演出 ---
$wsShow.Cells.Item(2, 6).Value = 28

# --- Sheet 4: 全部类型 ---
$wsAll.Cells.Item(2, 6).Value = 4738
$wsAll.Cells.Item(3, 6).Value = 1890
$wsAll.Cells.Item(6, 6).Value = 3192
$wsAll.Cells.Item(7, 6).Value = 586
$wsAll.Cells.Item(8, 6).Value = 606
$wsAll.Cells.Item(9, 6).Value = 293
$wsAll.Cells.Item(10, 6).Value = 659
$wsAll.Cells.Item(11, 6).Value = 555
$wsAll.Cells.Item(12, 6).Value = 561
$wsAll.Cells.Item(13, 6).Value = 28
$wsAll.Cells.Item(17, 6).Value = 1393
$wsAll.Cells.Item(19, 6).Value = 1655
$wsAll.Cells.Item(33, 6).Value = 4048
$wsAll.Cells.Item(35, 6).Value = 15
$wsAll.Cells.Item(39, 6).Value = 1716
$wsAll.Cells.Item(40, 6).Value = 61
$wsAll.Cells.Item(41, 6).Value = 1915
